$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.708.35'
$ws.Range("E2").Value = '  +1.80%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.869.48'
$ws.Range("E3").Value = '  +1.76%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9994'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.38'
$ws.Range("E5").Value = '  +1.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9994'
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4782'
$ws.Range("E7").Value = '  +1.99%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2826'
$ws.Range("E8").Value = '  +4.59%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06513'
$ws.Range("E9").Value = '  +3.95%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.82'
$ws.Range("E10").Value = '  +17.53%  '
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07493'
$ws.Range("E11").Value = '  +1.01%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.828.84'
$ws.Range("E12").Value = '  -0.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '94.53'
$ws.Range("E13").Value = '  +13.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.079'
$ws.Range("E14").Value = '  +3.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6492'
$ws.Range("E15").Value = '  +4.90%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '295.52'
$ws.Range("E16").Value = '  +29.97%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.660.13'
$ws.Range("E17").Value = '  +1.92%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9986'
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.99'
$ws.Range("E19").Value = '  +5.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007479'
$ws.Range("E20").Value = '  +2.83%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.107.77'
$ws.Range("E21").Value = '  +1.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.171'
$ws.Range("E23").Value = '  +6.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.079'
$ws.Range("E24").Value = '  +4.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '168.74'
$ws.Range("E25").Value = '  +2.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.208'
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.41'
$ws.Range("E27").Value = '  +9.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.952'
$ws.Range("E28").Value = '  +4.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1044'
$ws.Range("E29").Value = '  +1.67%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.351'
$ws.Range("E30").Value = '  -1.42%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.112'
$ws.Range("E31").Value = '  +0.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.944'
$ws.Range("E32").Value = '  +4.18%  '
$ws.Range("E33").Value = '  +3.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.171'
$ws.Range("E34").Value = '  +3.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7176'
$ws.Range("E35").Value = '  +1.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.703'
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01927'
$ws.Range("E37").Value = '  +3.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.700'
$ws.Range("E38").Value = '  +1.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.033'
$ws.Range("E39").Value = '  +5.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8898'
$ws.Range("E40").Value = '  -0.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '107.16'
$ws.Range("E41").Value = '  +2.73%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9993'
$ws.Range("E42").Value = '  -0.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4187'
$ws.Range("E43").Value = '  +4.72%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.560'
$ws.Range("E44").Value = '  +0.56%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.319'
$ws.Range("E45").Value = '  +5.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.36'
$ws.Range("E46").Value = '  +8.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1225'
$ws.Range("E47").Value = '  +2.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '34.56'
$ws.Range("E48").Value = '  +6.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.786'
$ws.Range("E49").Value = '  +3.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.390'
$ws.Range("E50").Value = '  +2.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05547'
$ws.Range("E51").Value = '  +0.76%  '
